$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: make the (already-default) width explicit, as in the target file ---
$ws.Columns("C").ColumnWidth = $ws.Columns("C").ColumnWidth

# --- Clear the "still needs review" breakdown values in the Round Robin table ---
# Row 30 (P3''): only E30 and G30 got cleared, D30/F30 keep their values
$ws.Range("E30").ClearContents()
$ws.Range("G30").ClearContents()

# Rows 31-33 (P3', p3, p3 follow-ups): the whole D:G block got cleared
$ws.Range("D31:G33").ClearContents()

# --- View state: scroll down toward the Round Robin table, zoom in, and land the
#     selection on E30 (the cell the author was re-checking) ---
$ws.Range("E30").Select()
$win = $excel.ActiveWindow
$win.Zoom = 160
